# Auto-generated: refresh the cryptos price/volume table to match the
# Sun Oct  8 16:16:45 UTC 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.892.29"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.631.77"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'211.42"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'23.48"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "1.864.05"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "1.633.51"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").Value = "'0.563"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").Value = "'65.46"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "27.897.40"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "'229.88"
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").Value = "'7.69"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").Value = "0.0₃0720"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").Value = "'10.07"
$ws.Range("E23").Value = "  -3.78%  "
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("D25").Value = "'155.00"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'15.51"
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").Value = "'0.0482"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("D34").Value = "1.392.83"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").Value = "'1.59"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  +10.17%  "
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").Value = "'0.850"
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'1.01"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").Value = "'65.76"
$ws.Range("E44").Value = "  -1.98%  "
$ws.Range("D45").Value = "'5.44"
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("D46").Value = "1.774.16"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("E47").Value = "  -2.74%  "
$ws.Range("D48").Value = "'88.65"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").Value = "'7.65"
$ws.Range("E51").Value = "  +1.20%  "
